# "User can set derived params"
# Replace the two GM_300 / GM_600 sample experiments (exp_ID 1-4, rows 2-5)
# with the three WW_Samanta_DOM experiments (exp_ID 48-50, rows 2-4), and
# drop the stray formatted placeholder row (row 10).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Remove the stray placeholder row (E10, gray Consolas font) and the
#     now-unused 4th data row. Delete the higher row index first so the
#     lower one's index doesn't shift before we get to it.
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(5).Delete()

# --- Helper: write a plain number into a cell that is formatted as Text
#     (numFmtId 49, "@") without Excel coercing the value into a string.
#     NumberFormat = "general" round-trips onto the existing built-in
#     General format (no new numFmt/cellXfs entries get minted), unlike
#     "General" which the host treats as a brand-new custom format.
function Set-NumericValue($range, $value) {
    $fmt = $range.NumberFormat
    $range.NumberFormat = "general"
    $range.Value = $value
    $range.NumberFormat = $fmt
}

# --- Row 2: exp_ID 48, WW_Samanta_DOM_2006
$ws.Range("A2").Value = 48
$ws.Range("B2").Value = "sim_WW_Samanta_DOM_2006.json"
$ws.Range("C2").Value = "crop_WW_Samanta_DOM_2006.json"
$ws.Range("D2").Value = "dom.json"
$ws.Range("E2").Value = "dom.csv"
$ws.Range("F2").Value = "wheat_CZ.json"
$ws.Range("G2").Value = "winter wheat_CZ.json"
Set-NumericValue $ws.Range("H2") 0
$ws.Range("I2").Value = "WW"

# --- Row 3: exp_ID 49, WW_Samanta_DOM_2007
$ws.Range("A3").Value = 49
$ws.Range("B3").Value = "sim_WW_Samanta_DOM_2007.json"
$ws.Range("C3").Value = "crop_WW_Samanta_DOM_2007.json"
$ws.Range("D3").Value = "dom.json"
$ws.Range("E3").Value = "dom.csv"
$ws.Range("F3").Value = "wheat_CZ.json"
$ws.Range("G3").Value = "winter wheat_CZ.json"
Set-NumericValue $ws.Range("H3") 0
$ws.Range("I3").Value = "WW"

# --- Row 4: exp_ID 50, WW_Samanta_DOM_2008 (previously a "1-2-3" rotation
#     row with text-typed H4 -- now a plain derived-params row like 2 & 3)
$ws.Range("A4").Value = 50
$ws.Range("B4").Value = "sim_WW_Samanta_DOM_2008.json"
$ws.Range("C4").Value = "crop_WW_Samanta_DOM_2008.json"
$ws.Range("D4").Value = "dom.json"
$ws.Range("E4").Value = "dom.csv"
$ws.Range("F4").Value = "wheat_CZ.json"
$ws.Range("G4").Value = "winter wheat_CZ.json"
Set-NumericValue $ws.Range("H4") 0
$ws.Range("I4").Value = "WW"

# --- Column widths: the sim/crop file names are now much longer, so widen
#     B & C to fit.
$ws.Columns.Item(2).ColumnWidth = 32.14
$ws.Columns.Item(3).ColumnWidth = 32.86

# --- Reset the view: move the selection off the old multi-cell block onto
#     a single cell below the data (also drops the stale topLeftCell scroll
#     position from when column E was scrolled into view).
$ws.Range("E13").Select()
